# typo in home page
# Applies the DeviceCatalog.pptx "home page" slide edits:
#  1. Resize/reposition + retext the "A medical device" -> "A model of medical device" label.
#  2. Restyle the "parentDevice (is part of)" label: strike-through + red for the
#     existing text, and append a new " hasPart" (green) run.
#  3. Reposition four small shapes around the "ChargeItemDefinition" / "Billing
#     extension" cluster (connector, rectangle, two "0..*" labels) to match the
#     relocated box.
#
# NOTE ON COORDINATES: PowerPoint's Shape.Left/Top/Width/Height are expressed in
# points (EMU / 12700) and are backed by a single-precision (float32) COM
# property. A naive `emu / 12700.0` literal can therefore land a hair on the
# wrong side of the float32 bucket boundary and save back 1 EMU off target, so
# the literals below were solved to land solidly inside the bucket that round-
# trips to the exact target EMU.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# 1) "A medical device" textbox (shape id=4, "ZoneTexte 3")
# ---------------------------------------------------------------------------
$shpDevice = $s.Shapes.Item(1)
$shpDevice.Left = 413.5136413574219    # 5251623 EMU
$shpDevice.Top = 194.96331787109375    # 2476034 EMU
$shpDevice.Width = 126.3065414428711   # 1604093 EMU
$shpDevice.Height = 50.892208099365234 # 646331 EMU (unchanged, set for completeness)
$shpDevice.TextFrame.TextRange.Text = "A model of medical device"

# ---------------------------------------------------------------------------
# 2) "parentDevice (is part of)" textbox (shape id=64, "ZoneTexte 63")
# ---------------------------------------------------------------------------
$shpParent = $s.Shapes.Item(7)
$trParent = $shpParent.TextFrame.TextRange

# Append " " then "hasPart" as two new runs *before* touching runs 1/2 so that
# they inherit run 2's pristine (italic-only) formatting rather than the
# strike/red styling we are about to add.
$null = $trParent.InsertAfter(" ")
$null = $trParent.InsertAfter("hasPart")

# New run 3: " " -> plain italic (inherited), no color override.
$run3 = $trParent.Runs(3)
$run3.Font.Italic = $true

# New run 4: "hasPart" -> italic, green.
$run4 = $trParent.Runs(4)
$run4.Font.Italic = $true
$run4.Font.Color.RGB = 5287936    # RGB(0,176,80) -> srgbClr 00B050

# Run 1: "parentDevice " -> strike-through, red
$run1 = $trParent.Runs(1)
$run1.Font.Strike = $true
$run1.Font.Color.RGB = 255        # RGB(255,0,0) -> srgbClr FF0000

# Run 2: "(is part of)" -> strike-through, red (keeps its existing italics)
$run2 = $trParent.Runs(2)
$run2.Font.Strike = $true
$run2.Font.Color.RGB = 255        # RGB(255,0,0) -> srgbClr FF0000

# ---------------------------------------------------------------------------
# 3) "0..*" textbox near parentDevice (shape id=44, "ZoneTexte 43")
# ---------------------------------------------------------------------------
$shp44 = $s.Shapes.Item(25)
$shp44.Left = 541.991455078125   # 6883291 EMU
$shp44.Top = 203.2296905517578   # 2581017 EMU

# ---------------------------------------------------------------------------
# 4) Bent connector above "ChargeItemDefinition" (shape id=49)
# ---------------------------------------------------------------------------
$shp49 = $s.Shapes.Item(34)
$shp49.Left = 185.14276123046875   # 2351313 EMU
$shp49.Width = 180.48606872558594  # 2292173 EMU
$shp49.Height = 34.66614532470703  # 440260 EMU

# ---------------------------------------------------------------------------
# 5) "ChargeItemDefinition" rectangle (shape id=50, "Rectangle 49")
# ---------------------------------------------------------------------------
$shp50 = $s.Shapes.Item(35)
$shp50.Left = 89.06764221191406    # 1131159 EMU
$shp50.Top = 356.5386047363281     # 4528040 EMU
$shp50.Width = 192.15000915527344  # 2440305 EMU

# ---------------------------------------------------------------------------
# 6) "0..*" textbox near ChargeItemDefinition (shape id=51, "ZoneTexte 50")
# ---------------------------------------------------------------------------
$shp51 = $s.Shapes.Item(36)
$shp51.Left = 181.81260681152344  # 2309020 EMU
$shp51.Top = 329.1330871582031    # 4179990 EMU
